$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking strings (e.g. "579.35") that must
# remain plain text (matching the original workbook inline-string cells).
# Force text format before assigning so Excel does not coerce them to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.862.59"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.608.07"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.35"
$ws.Range("E5").Value = "  +4.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.83"
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.633.63"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("E10").Value = "  -2.90%  "
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("E12").Value = "  -5.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.371"
$ws.Range("E13").Value = "  +5.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.078.03"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.827.80"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.624.16"
$ws.Range("E18").Value = "  +1.61%  "
$ws.Range("E19").Value = "  +9.49%  "
$ws.Range("E20").Value = "  +3.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.50"
$ws.Range("E21").Value = "  +3.82%  "
$ws.Range("E22").Value = "  +7.24%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.519"
$ws.Range("E24").Value = "  +8.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.21"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.93"
$ws.Range("E28").Value = "  +7.57%  "
$ws.Range("E29").Value = "  +3.55%  "
$ws.Range("E30").Value = "  +9.17%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.36"
$ws.Range("E32").Value = "  +2.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "162.65"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("E34").Value = "  +2.75%  "
$ws.Range("E35").Value = "  +13.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.28"
$ws.Range("E36").Value = "  +4.90%  "
$ws.Range("E37").Value = "  +6.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.62"
$ws.Range("E38").Value = "  +9.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.98"
$ws.Range("E40").Value = "  +5.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "309.55"
$ws.Range("E41").Value = "  +7.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.849"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "134.22"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.51"
$ws.Range("E44").Value = "  +10.13%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.06"
$ws.Range("E45").Value = "  +12.24%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.96"
$ws.Range("E46").Value = "  +5.88%  "
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("E48").Value = "  +2.88%  "
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0551"
$ws.Range("E50").Value = "  +3.91%  "
$ws.Range("E51").Value = "  +3.95%  "
